$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.89"
$ws.Range("E2").Value = "'1.69%"

$ws.Range("D3").Value = "'37.72"
$ws.Range("E3").Value = "'0.84%"

$ws.Range("D4").Value = "'5.123"
$ws.Range("E4").Value = "'1.08%"

$ws.Range("D5").Value = "'0.07881"
$ws.Range("E5").Value = "'2.00%"

$ws.Range("D6").Value = "'4.421"
$ws.Range("E6").Value = "'1.84%"

$ws.Range("D7").Value = "'8.268"
$ws.Range("E7").Value = "'1.02%"

$ws.Range("D8").Value = "'1.906"
$ws.Range("E8").Value = "'0.34%"

$ws.Range("D9").Value = "'2.837"
$ws.Range("E9").Value = "'-10.61%"

$ws.Range("D10").Value = "'0.9213"
$ws.Range("E10").Value = "'0.15%"

$ws.Range("D11").Value = "'0.1180"
$ws.Range("E11").Value = "'-4.50%"

$ws.Range("D12").Value = "'0.1926"
$ws.Range("E12").Value = "'3.12%"

$ws.Range("D13").Value = "'0.09019"
$ws.Range("E13").Value = "'2.20%"

$ws.Range("D14").Value = "'0.03316"
$ws.Range("E14").Value = "'-3.02%"

$ws.Range("D15").Value = "'0.09605"
$ws.Range("E15").Value = "'-0.99%"

$ws.Range("E16").Value = "'1.33%"

$ws.Range("D17").Value = "'0.005991"
$ws.Range("E17").Value = "'-0.29%"

$ws.Range("E18").Value = "'-0.44%"

$ws.Range("D19").Value = "'0.3443"
$ws.Range("E19").Value = "'1.00%"

$ws.Range("D20").Value = "'5.237"
$ws.Range("E20").Value = "'4.44%"

$ws.Range("D21").Value = "'0.1285"
$ws.Range("E21").Value = "'1.23%"

$ws.Range("D23").Value = "'0.04356"
$ws.Range("E23").Value = "'0.58%"

$ws.Range("D24").Value = "'0.001248"
$ws.Range("E24").Value = "'3.05%"

$ws.Range("D25").Value = "'0.004663"
$ws.Range("E25").Value = "'10.30%"

$ws.Range("D26").Value = "'0.0001359"
$ws.Range("E26").Value = "'0.66%"

$ws.Range("D27").Value = "'0.0003988"

$ws.Range("D39").Value = "'0.02247"
$ws.Range("E39").Value = "'3.44%"

$ws.Range("D40").Value = "'0.05075"
$ws.Range("E40").Value = "'3.66%"

$ws.Range("D41").Value = "'0.007460"
$ws.Range("E41").Value = "'-2.57%"

$ws.Range("D42").Value = "'0.009039"
$ws.Range("E42").Value = "'-8.98%"

$ws.Range("D43").Value = "'0.1354"
$ws.Range("E43").Value = "'1.28%"

$ws.Range("D44").Value = "'0.001948"
$ws.Range("E44").Value = "'-2.29%"

$ws.Range("D45").Value = "'0.008609"
$ws.Range("E45").Value = "'-12.59%"

$ws.Range("D46").Value = "'0.00006587"
$ws.Range("E46").Value = "'0.54%"

$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.07%"

$ws.Range("D48").Value = "'0.003306"
$ws.Range("E48").Value = "'10.15%"

$ws.Range("E49").Value = "'-23.14%"

$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.07%"

$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.07%"
